# словарьТЗ.xlsx - apply commit "parallel operation and correct file names are configured"
#
# Summary of the change (from the OOXML diff):
#  - Sheet "Титульный лист": two new rows are inserted at the top of the
#    data block: "year" (moved up from its old spot) and a brand-new
#    "second_name" entry. Everything that used to start at row 4
#    (university, faculty, ... group_number, year) shifts down by two rows,
#    and the now-duplicated "year" row further down is cleared out.
#  - The "akad_post" helper text had a stray embedded line break that is
#    removed (becomes one continuous line).
#  - The active sheet on open changes from "Пояснительная записка" to
#    "Титульный лист", and the selection / scroll position of a couple of
#    sheets changes accordingly.

$wb = $excel.ActiveWorkbook

$wsTitle = $wb.Worksheets.Item("Титульный лист")
$wsExpl  = $wb.Worksheets.Item("Пояснительная записка")

# ---------------------------------------------------------------------
# 1. Make room for the two new rows (university ... year all shift +2)
# ---------------------------------------------------------------------
$wsTitle.Range("A4:A5").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2. New row 4: "year" (promoted from the bottom of the list)
# ---------------------------------------------------------------------
$wsTitle.Range("A4").Value = "year"
$wsTitle.Range("B4").Value = "Введите год"
$wsTitle.Range("C4").Formula = "=""'""&A4&""': '',"""
$wsTitle.Range("D4").Formula = "=""('""&B4&"":', '""&A4&""'),"""

# ---------------------------------------------------------------------
# 3. New row 5: brand new "second_name" field
# ---------------------------------------------------------------------
$wsTitle.Range("A5").Value = "second_name"
$wsTitle.Range("B5").Value = 'Введите Фамилию на английском языке, например "Shapovalov"'
$wsTitle.Range("C5").Formula = "=""'""&A5&""': '',"""
$wsTitle.Range("D5").Formula = "=""('""&B5&"":', '""&A5&""'),"""

# ---------------------------------------------------------------------
# 4. The old "year" row (now sitting two rows further down, at row 18)
#    is now redundant since "year" lives at row 4 - clear it out.
# ---------------------------------------------------------------------
$wsTitle.Range("A18:D18").ClearContents()

# ---------------------------------------------------------------------
# 5. Fix the "akad_post" helper text: drop the stray line break so the
#    whole hint reads as a single line (row 11 shifted to row 13).
# ---------------------------------------------------------------------
$wsTitle.Range("B13").Value = 'Введите должность руководителя, утвердившего документ, например "Академический руководитель образовательной программы «Программная инженерия», кандидат технических наук"'

# ---------------------------------------------------------------------
# 6. View / active-sheet bookkeeping: "Титульный лист" becomes the sheet
#    that is active when the workbook is reopened.
# ---------------------------------------------------------------------
$wsTitle.Select()
$wsTitle.Range("D15").Select()

$wsExpl.Range("B3").Select()

$wb.Windows.Item(1).ActiveSheet
